$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data Collection")

$data = @(
    ,@("@oracle", 43051, 129, 0.146969417318254, 0.0682694889090238, 43051.686574074076, 28625, 53153, 1126288, 2)
    ,@("@oracle", 43050, 279, 0.0825575131757927, 0.0333644060290565, 43050.74318287037, 79235, 147363, 1164567, 13)
    ,@("@oracle", 43049, 176, 0.271843864784205, 0.141242662009707, 43049.7084375, 2328, 4873, 903368, 9)
    ,@("@microsoft", 43051, 436, 0.343027030627374, 0.126860207994095, 43051.73405092592, 18030, 47139, 2382266, 11)
    ,@("@microsoft", 43050, 589, 0.353411599012193, 0.161388989613038, 43050.70775462963, 26334, 91877, 2926112, 27)
    ,@("@microsoft", 43049, 1067, 0.31842923759395, 0.123030835252601, 43049.72518518518, 54681, 159766, 16890622, 54)
    ,@("@microsoft", 43048, 1350, 0.346868158213528, 0.12900209803057, 43048.73851851852, 132180, 472977, 12224187, 49)
    ,@("@microsoft", 43047, 2589, 0.296453870096689, 0.133648923567509, 43047.727800925924, 886006, 5678150, 10942134, 69)
    ,@("@microsoft", 43046, 1884, 0.335308670005193, 0.120278928433129, 43046.73800925926, 67370, 455760, 15655615, 62)
    ,@("@microsoft", 43045, 1246, 0.337933145714044, 0.142007954323112, 43045.75975694445, 55237, 305135, 24069567, 81)
    ,@("@google", 43051, 1154, 0.266585324337707, 0.0896151476440145, 43051.73269675926, 621286, 1879527, 5276634, 44)
    ,@("@google", 43050, 1741, 0.292561855805393, 0.149620778610377, 43050.73390046296, 634177, 1717773, 18651373, 51)
    ,@("@google", 43049, 2410, 0.320259056756979, 0.204810433452339, 43049.70853009259, 728970, 1882366, 14133192, 99)
    ,@("@google", 43048, 2731, 0.279319951327, 0.118267685992971, 43048.74686342593, 528803, 1260719, 22383865, 107)
    ,@("@google", 43047, 2420, 0.357511044663935, 0.190866551140214, 43047.74260416667, 741334, 1855603, 17160497, 81)
    ,@("@google", 43046, 2662, 0.367301685952139, 0.181640862778576, 43046.75262731482, 889832, 2736449, 34626861, 119)
    ,@("@google", 43045, 3190, 0.421112982150205, 0.176021183929344, 43045.75341435185, 1322624, 3717511, 15641992, 87)
    ,@("@intel", 43051, 159, 0.217540862352183, 0.0894082829076933, 43051.74265046296, 16500, 68519, 913879, 3)
    ,@("@intel", 43050, 224, 0.284730081426509, 0.102158426580859, 43050.72318287037, 13316, 61867, 1100498, 6)
    ,@("@intel", 43049, 330, 0.316131440669319, 0.153106849519728, 43049.739583333336, 19854, 73869, 4077966, 14)
    ,@("@intel", 43048, 459, 0.4009750778351, 0.164504782038142, 43048.753703703704, 33208, 125272, 11336102, 45)
    ,@("@intel", 43047, 504, 0.474238594848714, 0.242112808928284, 43047.76008101852, 61916, 90765, 10820240, 36)
    ,@("@intel", 43046, 937, 0.488971349414251, 0.33080688151526, 43046.781273148146, 178084, 285854, 9479257, 69)
    ,@("@intel", 43045, 484, 0.298151428939135, 0.131109791314021, 43045.75096064815, 29272, 146553, 5138084, 30)
    ,@("@xerox", 43051, 15, 0.279132996632996, 0.02510101010101, 43051.68740740741, 318, 724, 511258, 1)
    ,@("@xerox", 43050, 21, 0.25883065347351, 0.0850065278636707, 43050.72414351852, 420, 862, 54216, 0)
    ,@("@xerox", 43049, 72, 0.272164351851851, 0.106828703703703, 43049.75048611111, 540, 1323, 3708280, 2)
    ,@("@xerox", 43048, 67, 0.335526910900045, 0.125075380672395, 43048.77888888889, 482, 559, 577345, 7)
    ,@("@xerox", 43047, 64, 0.274767203282828, 0.186375473484848, 43047.777916666666, 504, 320, 55167, 0)
    ,@("@xerox", 43046, 76, 0.33938230994152, 0.262171052631578, 43046.809328703705, 860, 587, 56512, 0)
    ,@("@xerox", 43045, 42, 0.182494588744588, 0.10340909090909, 43045.73509259259, 398, 598, 195108, 2)
    ,@("@cisco", 43051, 112, 0.239983057616986, 0.134898208810039, 43051.73981481481, 2032, 3489, 233753, 0)
    ,@("@cisco", 43050, 114, 0.304128271014235, -0.0143225013290802, 43050.76613425926, 3622, 5838, 2978744, 3)
    ,@("@cisco", 43049, 318, 0.30267656347845, 0.0221117168994527, 43049.74395833333, 5095, 5580, 4078743, 17)
    ,@("@cisco", 43048, 295, 0.338278254701983, 0.20487648924169, 43048.718877314815, 1423, 4078, 1407384, 11)
    ,@("@cisco", 43047, 365, 0.294350639467078, 0.14133754950022, 43047.739583333336, 3305, 7341, 1959983, 19)
    ,@("@cisco", 43046, 464, 0.259819115570731, 0.15519593448338, 43046.74119212963, 4224, 5594, 1855143, 22)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $ws.Cells.Item($r, 9).Value = $row[8]
    $ws.Cells.Item($r, 10).Value = $row[9]
}

$ws.Range("E34").Select()

Write-Output "done"